# BAU Components Energy Use - South Korea files update
# 1) Fix typo in unit-conversion label on the "About" sheet: "1 TOE to mBTU" -> "1 TOE to mmBTU"
# 2) Fix the mis-spelled "Snopshot2040" sheet name -> "Snapshot2040"
# 3) Restore view/selection state to match the author's last-saved position

$wb = $excel.ActiveWorkbook

# --- 1) Correct the TOE-to-mmBTU conversion label on the About sheet ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A30").Value = "1 TOE to mmBTU"

# --- 2) Rename the mis-spelled sheet ---
$wsSnapshot = $wb.Worksheets.Item("Snopshot2040")
$wsSnapshot.Name = "Snapshot2040"

# --- 3) Restore sheet selections / active sheet ---
$wsAbout.Activate()
$wsAbout.Range("A31").Select()

$wsForecast = $wb.Worksheets.Item("Forecast")
$wsForecast.Activate()
$wsForecast.Range("B10").Select()

# Snapshot2040 is the last-active / selected tab in the saved workbook
$wsSnapshot.Activate()
$wsSnapshot.Range("A43").Select()
